$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 10 cells to row 11 to replicate style (A column bold/border/center, E column blank inline string)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# Set cell values for the new row 11 (2021 data)
$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 1843.93
$ws.Range("C11").Value = 351.76
$ws.Range("D11").Value = 80.95999999999999
$ws.Range("F11").Value = 912.0599999999999
$ws.Range("G11").Value = 1953
$ws.Range("H11").Value = 402.35
$ws.Range("I11").Value = 1617.48
$ws.Range("J11").Value = 121.98
$ws.Range("K11").Value = 33431.63
$ws.Range("L11").Value = 174.64
$ws.Range("M11").Value = 51.9
$ws.Range("N11").Value = 15.41
$ws.Range("O11").Value = 720.3200000000001
$ws.Range("P11").Value = 1076.13
$ws.Range("Q11").Value = 54.44
$ws.Range("R11").Value = 106.12
$ws.Range("S11").Value = 826.65
$ws.Range("T11").Value = 14.51
$ws.Range("U11").Value = 2792
$ws.Range("V11").Value = 208.48
$ws.Range("W11").Value = 350.75
$ws.Range("X11").Value = 40.65
$ws.Range("Y11").Value = 28.2
$ws.Range("Z11").Value = 3371.96
$ws.Range("AA11").Value = 227.42
$ws.Range("AB11").Value = 1318.92
$ws.Range("AC11").Value = 94.65000000000001
$ws.Range("AD11").Value = 773.0700000000001
$ws.Range("AE11").Value = 575.13
$ws.Range("AF11").Value = 4839.14
$ws.Range("AG11").Value = 1853.05
$ws.Range("AH11").Value = 383.9
$ws.Range("AI11").Value = 662.29
$ws.Range("AJ11").Value = 15.36
$ws.Range("AK11").Value = 1137.13
$ws.Range("AL11").Value = 661.47
$ws.Range("AM11").Value = 1123.86
$ws.Range("AN11").Value = 27
$ws.Range("AO11").Value = 445.4
$ws.Range("AP11").Value = 2096.75
$ws.Range("AQ11").Value = 81.25
